$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Repayment Schedule")

# Insert a new blank column before column N (shifts N:P -> O:Q)
$ws.Columns("N:N").Insert()

# Make this sheet the active tab with a new selection, matching the authored edit
$ws.Select()
$ws.Range("R7").Select()

# The previously-active sheet (NewLoanInput) should no longer be the selected tab
$ws1 = $wb.Worksheets.Item("NewLoanInput")
$ws1.Range("B2").Select()

$ws.Select()
$ws.Range("R7").Select()
